$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 11).Value = 1755
$ws.Cells.Item(3, 11).Value = 1673
$ws.Cells.Item(4, 10).Value = 1801
$ws.Cells.Item(4, 11).Value = 357
$ws.Cells.Item(5, 11).Value = 113
$ws.Cells.Item(6, 11).Value = 2156
$ws.Cells.Item(7, 10).Value = 29269
$ws.Cells.Item(7, 11).Value = 6054

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(2, 11).Value = 22
$ws.Cells.Item(3, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 95

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 11).Value = 118
$ws.Cells.Item(3, 11).Value = 116
$ws.Cells.Item(5, 11).Value = 10
$ws.Cells.Item(6, 11).Value = 135
$ws.Cells.Item(7, 11).Value = 398

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(3, 11).Value = 93
$ws.Cells.Item(7, 11).Value = 246

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(6, 11).Value = 74
$ws.Cells.Item(7, 11).Value = 207

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(3, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 151

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(2, 11).Value = 36
$ws.Cells.Item(3, 11).Value = 37
$ws.Cells.Item(6, 11).Value = 35
$ws.Cells.Item(7, 11).Value = 114

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 11).Value = 46
$ws.Cells.Item(4, 11).Value = 26
$ws.Cells.Item(7, 11).Value = 171
$ws.Cells.Item(8, 11).Value = 398
$ws.Cells.Item(10, 11).Value = 37
$ws.Cells.Item(11, 11).Value = 126
$ws.Cells.Item(16, 11).Value = 12
$ws.Cells.Item(18, 11).Value = 45
$ws.Cells.Item(19, 11).Value = 168
$ws.Cells.Item(20, 11).Value = 134
$ws.Cells.Item(22, 11).Value = 14
$ws.Cells.Item(25, 11).Value = 31
$ws.Cells.Item(29, 11).Value = 290
$ws.Cells.Item(31, 11).Value = 68
$ws.Cells.Item(33, 11).Value = 246
$ws.Cells.Item(36, 11).Value = 70
$ws.Cells.Item(37, 11).Value = 207
$ws.Cells.Item(41, 11).Value = 61
$ws.Cells.Item(44, 11).Value = 56
$ws.Cells.Item(48, 11).Value = 66
$ws.Cells.Item(51, 11).Value = 70
$ws.Cells.Item(52, 11).Value = 164
$ws.Cells.Item(53, 11).Value = 95
$ws.Cells.Item(63, 10).Value = 91
$ws.Cells.Item(63, 11).Value = 22
$ws.Cells.Item(65, 11).Value = 151
$ws.Cells.Item(67, 11).Value = 233
$ws.Cells.Item(73, 11).Value = 60
$ws.Cells.Item(77, 11).Value = 40
$ws.Cells.Item(78, 11).Value = 77
$ws.Cells.Item(79, 11).Value = 161
$ws.Cells.Item(80, 11).Value = 22
$ws.Cells.Item(85, 11).Value = 301
$ws.Cells.Item(88, 11).Value = 74
$ws.Cells.Item(90, 11).Value = 52
$ws.Cells.Item(92, 11).Value = 27
$ws.Cells.Item(96, 11).Value = 85
$ws.Cells.Item(97, 11).Value = 53
$ws.Cells.Item(99, 11).Value = 114
$ws.Cells.Item(101, 10).Value = 29269
$ws.Cells.Item(101, 11).Value = 6054

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(7, 11).Value = 68

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(6, 11).Value = 78
$ws.Cells.Item(7, 11).Value = 233

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 11).Value = 79
$ws.Cells.Item(3, 11).Value = 98
$ws.Cells.Item(6, 11).Value = 93
$ws.Cells.Item(7, 11).Value = 290

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(2, 11).Value = 15
$ws.Cells.Item(6, 11).Value = 30
$ws.Cells.Item(7, 11).Value = 66

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 11).Value = 53
$ws.Cells.Item(3, 11).Value = 50
$ws.Cells.Item(6, 11).Value = 54
$ws.Cells.Item(7, 11).Value = 168

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(2, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 56

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(4, 11).Value = 5
$ws.Cells.Item(6, 11).Value = 50

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Cells.Item(2, 11).Value = 20
$ws.Cells.Item(6, 11).Value = 27
$ws.Cells.Item(7, 11).Value = 61

$ws = $wb.Worksheets.Item("Avondale")
$ws.Cells.Item(3, 11).Value = 5
$ws.Cells.Item(7, 11).Value = 37

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(3, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 77

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(3, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 85

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(4, 11).Value = 12
$ws.Cells.Item(7, 11).Value = 161

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(6, 11).Value = 49
$ws.Cells.Item(7, 11).Value = 134

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(3, 11).Value = 14
$ws.Cells.Item(7, 11).Value = 45

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(3, 11).Value = 25
$ws.Cells.Item(7, 11).Value = 70

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(5, 11).Value = 7
$ws.Cells.Item(6, 11).Value = 43
$ws.Cells.Item(7, 11).Value = 171

$ws = $wb.Worksheets.Item("East Side")
$ws.Cells.Item(3, 11).Value = 12
$ws.Cells.Item(7, 11).Value = 31

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(3, 11).Value = 32
$ws.Cells.Item(7, 11).Value = 126

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(2, 11).Value = 16
$ws.Cells.Item(7, 11).Value = 60

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(2, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 46

$ws = $wb.Worksheets.Item("West Town")
$ws.Cells.Item(3, 11).Value = 7
$ws.Cells.Item(6, 11).Value = 35
$ws.Cells.Item(7, 11).Value = 53

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Cells.Item(6, 11).Value = 16
$ws.Cells.Item(7, 11).Value = 27

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(3, 11).Value = 14
$ws.Cells.Item(7, 11).Value = 74

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(3, 11).Value = 12
$ws.Cells.Item(7, 11).Value = 52

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(2, 11).Value = 17
$ws.Cells.Item(3, 11).Value = 20
$ws.Cells.Item(7, 11).Value = 70

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(3, 11).Value = 98
$ws.Cells.Item(5, 11).Value = 5
$ws.Cells.Item(7, 11).Value = 301

$ws = $wb.Worksheets.Item("Clearing")
$ws.Cells.Item(2, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 14

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(7, 11).Value = 40

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Cells.Item(2, 11).Value = 5
$ws.Cells.Item(7, 11).Value = 22

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(2, 11).Value = 37
$ws.Cells.Item(3, 11).Value = 36
$ws.Cells.Item(6, 11).Value = 77
$ws.Cells.Item(7, 11).Value = 164

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Cells.Item(3, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 26

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(7, 11).Value = 12
